$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): extend B1:O1 (values 0..13) to B1:Q1 (values 0..15)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header style (bold, centered, bordered) from O1 onto the new P1:Q1 cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-25: updated simulation results (columns B,C,D,F,G,I) plus new
# P (always 0) and Q (new remainder/tail value) columns; O now holds 0.
# Row 2
$ws.Range("B2").Value = 3.474798344480178
$ws.Range("C2").Value = 1.028811158724977
$ws.Range("D2").Value = 0.09244169505385003
$ws.Range("F2").Value = 0.215384746749443
$ws.Range("G2").Value = 0.0007878213449589255
$ws.Range("I2").Value = 0.004342034164140696
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.5037206014540914

# Row 3
$ws.Range("B3").Value = 3.028242116624824
$ws.Range("C3").Value = 0.9095836216172586
$ws.Range("D3").Value = 0.08144680489036915
$ws.Range("F3").Value = 0.2079164474233188
$ws.Range("G3").Value = 0.0007913191441614731
$ws.Range("I3").Value = 0.002627410753965886
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0.5170899731007097

# Row 4
$ws.Range("B4").Value = 2.753769637390974
$ws.Range("C4").Value = 0.8369175960381483
$ws.Range("D4").Value = 0.0747359746587648
$ws.Range("F4").Value = 0.2041134526760544
$ws.Range("G4").Value = 0.0007935369428712949
$ws.Range("I4").Value = 0.001813952300414812
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0.5279113412277923

# Row 5
$ws.Range("B5").Value = 2.641799028473429
$ws.Range("C5").Value = 0.8090124634292692
$ws.Range("D5").Value = 0.07213695982333945
$ws.Range("F5").Value = 0.2022530400963163
$ws.Range("G5").Value = 0.0007944626045199381
$ws.Range("I5").Value = 0.001596515654408215
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0.5315560935065662

# Row 6
$ws.Range("B6").Value = 2.623166834794176
$ws.Range("C6").Value = 0.8063278378857319
$ws.Range("D6").Value = 0.07186117016502891
$ws.Range("F6").Value = 0.2013499986950862
$ws.Range("G6").Value = 0.0007946221990319479
$ws.Range("I6").Value = 0.001642483165760211
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0.5305089285972784

# Row 7
$ws.Range("B7").Value = 2.75216823782398
$ws.Range("C7").Value = 0.8418601316625143
$ws.Range("D7").Value = 0.07512620537927717
$ws.Range("F7").Value = 0.2024359241411631
$ws.Range("G7").Value = 0.0007935620479241891
$ws.Range("I7").Value = 0.002010764238572271
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0.5233522562144231

# Row 8
$ws.Range("B8").Value = 3.320746749036402
$ws.Range("C8").Value = 0.9947037192501682
$ws.Range("D8").Value = 0.08921112445907653
$ws.Range("F8").Value = 0.2104624602156306
$ws.Range("G8").Value = 0.0007890280355672963
$ws.Range("I8").Value = 0.003933701769336828
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0.5017252544370336

# Row 9
$ws.Range("B9").Value = 4.434762411754207
$ws.Range("C9").Value = 1.289645937344005
$ws.Range("D9").Value = 0.1164251945498194
$ws.Range("F9").Value = 0.2359253096864506
$ws.Range("G9").Value = 0.0007806592506337806
$ws.Range("I9").Value = 0.009528435898541687
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0.4899511050140291

# Row 10
$ws.Range("B10").Value = 5.254271404277119
$ws.Range("C10").Value = 1.512968138794918
$ws.Range("D10").Value = 0.1369357146427319
$ws.Range("F10").Value = 0.2588012297296487
$ws.Range("G10").Value = 0.0007748235914224699
$ws.Range("I10").Value = 0.01534737945539533
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0.4953909904036493

# Row 11
$ws.Range("B11").Value = 5.627597724287341
$ws.Range("C11").Value = 1.62477588229757
$ws.Range("D11").Value = 0.1470796374983223
$ws.Range("F11").Value = 0.2678473146346505
$ws.Range("G11").Value = 0.0007722429752540227
$ws.Range("I11").Value = 0.01872997709034241
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0.4948416215512879

# Row 12
$ws.Range("B12").Value = 5.769190647156222
$ws.Range("C12").Value = 1.662452287175142
$ws.Range("D12").Value = 0.1505480697081367
$ws.Range("F12").Value = 0.2728827059642285
$ws.Range("G12").Value = 0.0007712671257639721
$ws.Range("I12").Value = 0.01991643805048859
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0.4991265633185975

# Row 13
$ws.Range("B13").Value = 5.738706857486704
$ws.Range("C13").Value = 1.65335173492474
$ws.Range("D13").Value = 0.1497223277782069
$ws.Range("F13").Value = 0.2720711800320714
$ws.Range("G13").Value = 0.0007714755825004936
$ws.Range("I13").Value = 0.01962236130179384
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0.4989409606163946

# Row 14
$ws.Range("B14").Value = 5.639251098238447
$ws.Range("C14").Value = 1.627456661801148
$ws.Range("D14").Value = 0.1473315415494625
$ws.Range("F14").Value = 0.2683773270195218
$ws.Range("G14").Value = 0.0007721619481533582
$ws.Range("I14").Value = 0.01881082380847499
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0.4955067809692224

# Row 15
$ws.Range("B15").Value = 5.578316350597504
$ws.Range("C15").Value = 1.613542739335628
$ws.Range("D15").Value = 0.146022594001991
$ws.Range("F15").Value = 0.2655860180574834
$ws.Range("G15").Value = 0.0007725861531889014
$ws.Range("I15").Value = 0.01839689215435136
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0.4919791961578568

# Row 16
$ws.Range("B16").Value = 5.229603591229932
$ws.Range("C16").Value = 1.521826652682591
$ws.Range("D16").Value = 0.1375640009747485
$ws.Range("F16").Value = 0.253538066662216
$ws.Range("G16").Value = 0.0007750172678775947
$ws.Range("I16").Value = 0.01570696817500039
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0.482698401191314

# Row 17
$ws.Range("B17").Value = 5.015883724479011
$ws.Range("C17").Value = 1.465770281291327
$ws.Range("D17").Value = 0.1323909025705774
$ws.Range("F17").Value = 0.2465075129474243
$ws.Range("G17").Value = 0.0007765242766214764
$ws.Range("I17").Value = 0.01417780346684516
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0.4781600032656854

# Row 18
$ws.Range("B18").Value = 4.893106207286223
$ws.Range("C18").Value = 1.42871326216482
$ws.Range("D18").Value = 0.1290308034419496
$ws.Range("F18").Value = 0.2440092956714182
$ws.Range("G18").Value = 0.0007773894687673949
$ws.Range("I18").Value = 0.01312846072559104
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0.4798235662855177

# Row 19
$ws.Range("B19").Value = 4.851479932017412
$ws.Range("C19").Value = 1.419946660042513
$ws.Range("D19").Value = 0.1281949237475857
$ws.Range("F19").Value = 0.2420897649791769
$ws.Range("G19").Value = 0.0007776891788000883
$ws.Range("I19").Value = 0.01294034349131579
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0.4774702580048995

# Row 20
$ws.Range("B20").Value = 5.038634438275665
$ws.Range("C20").Value = 1.471345795093384
$ws.Range("D20").Value = 0.1329103356176944
$ws.Range("F20").Value = 0.2473539082364624
$ws.Range("G20").Value = 0.0007763626748547517
$ws.Range("I20").Value = 0.01431928511958613
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0.4789034281026829

# Row 21
$ws.Range("B21").Value = 5.668362351303585
$ws.Range("C21").Value = 1.640548625708448
$ws.Range("D21").Value = 0.1484718562840612
$ws.Range("F21").Value = 0.2678756722076869
$ws.Range("G21").Value = 0.0007719676543110315
$ws.Range("I21").Value = 0.01925839013490371
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0.4922019874233428

# Row 22
$ws.Range("B22").Value = 6.08086761876217
$ws.Range("C22").Value = 1.744903541909025
$ws.Range("D22").Value = 0.1581426640339174
$ws.Range("F22").Value = 0.2845486363314151
$ws.Range("G22").Value = 0.0007691350043418969
$ws.Range("I22").Value = 0.0226357147787466
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0.5103524705731388

# Row 23
$ws.Range("B23").Value = 5.86073621393308
$ws.Range("C23").Value = 1.682516017764442
$ws.Range("D23").Value = 0.1524469523580478
$ws.Range("F23").Value = 0.2774363152505543
$ws.Range("G23").Value = 0.000770633497717991
$ws.Range("I23").Value = 0.0205555905296313
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0.5054626737246934

# Row 24
$ws.Range("B24").Value = 5.028530598248324
$ws.Range("C24").Value = 1.458486462537508
$ws.Range("D24").Value = 0.1318498588219654
$ws.Range("F24").Value = 0.249978908953473
$ws.Range("G24").Value = 0.0007764200885551076
$ws.Range("I24").Value = 0.01384744488613876
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0.4867696337544345

# Row 25
$ws.Range("B25").Value = 4.133163646007404
$ws.Range("C25").Value = 1.219142006372522
$ws.Range("D25").Value = 0.109807504906712
$ws.Range("F25").Value = 0.2252774970166698
$ws.Range("G25").Value = 0.0007828848980073834
$ws.Range("I25").Value = 0.00809595745542957
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0.4823595699843537
